# Auto-generated PowerShell Excel COM-interop script
# Applies numeric corrections to several Leve profit rows across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets (scheduled runner update).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 241.33333
$ws.Range("I11").Value = 241.33333
$ws.Range("K11").Value = 241.33333
$ws.Range("M11").Value = -101.33333

$ws.Range("H32").Value = 4699.0557
$ws.Range("J32").Value = 4398.923
$ws.Range("L32").Value = 4398.923
$ws.Range("N32").Value = -5050.923

$ws.Range("H33").Value = 15399.762
$ws.Range("I33").Value = 16941.842
$ws.Range("K33").Value = 16941.842
$ws.Range("M33").Value = -16712.842

$ws.Range("H113").Value = 142861260
$ws.Range("I113").Value = 166669140
$ws.Range("K113").Value = 166669140
$ws.Range("M113").Value = -166665886

$ws.Range("H116").Value = 3714.2856
$ws.Range("I116").Value = 3600
$ws.Range("J116").Value = 3800
$ws.Range("K116").Value = 3600
$ws.Range("L116").Value = 3800
$ws.Range("M116").Value = -158
$ws.Range("N116").Value = -10684

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 33450
$ws.Range("J55").Value = 33450
$ws.Range("L55").Value = 33450
$ws.Range("N55").Value = -34080

$ws.Range("H61").Value = 23816508
$ws.Range("I61").Value = 27783428
$ws.Range("K61").Value = 27783428
$ws.Range("M61").Value = -27783216

$ws.Range("H110").Value = 4641.1787
$ws.Range("I110").Value = 3826.8823
$ws.Range("K110").Value = 3826.8823
$ws.Range("M110").Value = -1781.8823

$ws.Range("H127").Value = 45000
$ws.Range("J127").Value = 45000
$ws.Range("L127").Value = 45000
$ws.Range("N127").Value = -54920

$ws.Range("H136").Value = 23816508
$ws.Range("I136").Value = 27783428
$ws.Range("K136").Value = 83350284
$ws.Range("M136").Value = -83347734

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 29800
$ws.Range("J35").Value = 29800
$ws.Range("L35").Value = 29800
$ws.Range("N35").Value = -30420

$ws.Range("H82").Value = 42250
$ws.Range("J82").Value = 42250
$ws.Range("L82").Value = 42250
$ws.Range("N82").Value = -43016

$ws.Range("H85").Value = 42250
$ws.Range("J85").Value = 42250
$ws.Range("L85").Value = 42250
$ws.Range("N85").Value = -44902

$ws.Range("H94").Value = 3731.7
$ws.Range("I94").Value = 3040
$ws.Range("K94").Value = 3040
$ws.Range("M94").Value = -2589

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 54489.8
$ws.Range("J28").Value = 54489.8
$ws.Range("L28").Value = 54489.8
$ws.Range("N28").Value = -54979.8

$ws.Range("H31").Value = 9099.700000000001
$ws.Range("I31").Value = 7285.4287
$ws.Range("K31").Value = 7285.4287
$ws.Range("M31").Value = -6990.4287

$ws.Range("H34").Value = 9099.700000000001
$ws.Range("I34").Value = 7285.4287
$ws.Range("K34").Value = 7285.4287
$ws.Range("M34").Value = -7083.4287

$ws.Range("H41").Value = 20266.111
$ws.Range("J41").Value = 19799.4
$ws.Range("L41").Value = 19799.4
$ws.Range("N41").Value = -20655.4

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H59").Value = 40561.75
$ws.Range("J59").Value = 40499
$ws.Range("L59").Value = 40499
$ws.Range("N59").Value = -42789

$ws.Range("H60").Value = 23500
$ws.Range("J60").Value = 24500
$ws.Range("L60").Value = 24500
$ws.Range("N60").Value = -25522

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H68").Value = 41092
$ws.Range("J68").Value = 41092
$ws.Range("L68").Value = 41092
$ws.Range("N68").Value = -42590

$ws.Range("H71").Value = 41092
$ws.Range("J71").Value = 41092
$ws.Range("L71").Value = 123276
$ws.Range("N71").Value = -130764

$ws.Range("H74").Value = 35558.332
$ws.Range("J74").Value = 35558.332
$ws.Range("L74").Value = 35558.332
$ws.Range("N74").Value = -37306.332

$ws.Range("H77").Value = 35558.332
$ws.Range("J77").Value = 35558.332
$ws.Range("L77").Value = 106674.996
$ws.Range("N77").Value = -115410.996

$ws.Range("H86").Value = 11062.333
$ws.Range("I86").Value = 11080.143
$ws.Range("J86").Value = 11000
$ws.Range("K86").Value = 11080.143
$ws.Range("L86").Value = 11000
$ws.Range("M86").Value = -9957.143
$ws.Range("N86").Value = -13246

$ws.Range("H89").Value = 11062.333
$ws.Range("I89").Value = 11080.143
$ws.Range("J89").Value = 11000
$ws.Range("K89").Value = 55400.715
$ws.Range("L89").Value = 55000
$ws.Range("M89").Value = -49784.715
$ws.Range("N89").Value = -66232

$ws.Range("H95").Value = 29999.666
$ws.Range("J95").Value = 29999.666
$ws.Range("L95").Value = 29999.666
$ws.Range("N95").Value = -35491.666

$ws.Range("H99").Value = 5533.3335
$ws.Range("I99").Value = 5750
$ws.Range("J99").Value = 5100
$ws.Range("K99").Value = 5750
$ws.Range("L99").Value = 5100
$ws.Range("M99").Value = -4252
$ws.Range("N99").Value = -8096

$ws.Range("H122").Value = 5164.875
$ws.Range("I122").Value = 4993.6665
$ws.Range("J122").Value = 5385
$ws.Range("K122").Value = 14980.9995
$ws.Range("L122").Value = 16155
$ws.Range("M122").Value = -12530.9995
$ws.Range("N122").Value = -21055

$ws.Range("H125").Value = 88662.5
$ws.Range("J125").Value = 88662.5
$ws.Range("L125").Value = 88662.5
$ws.Range("N125").Value = -93582.5

$ws.Range("H126").Value = 5533.3335
$ws.Range("I126").Value = 5750
$ws.Range("J126").Value = 5100
$ws.Range("K126").Value = 17250
$ws.Range("L126").Value = 15300
$ws.Range("M126").Value = -14780
$ws.Range("N126").Value = -20240

$ws.Range("H134").Value = 8304.23
$ws.Range("I134").Value = 7086.8184
$ws.Range("K134").Value = 21260.4552
$ws.Range("M134").Value = -18725.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2500
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2500
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 7500
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -8580

$ws.Range("H139").Value = 2604.5908
$ws.Range("I139").Value = 2604.5908
$ws.Range("K139").Value = 7813.7724
$ws.Range("M139").Value = -2673.7724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 8500
$ws.Range("J28").Value = 8500
$ws.Range("L28").Value = 8500
$ws.Range("N28").Value = -8884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3179.8
$ws.Range("I22").Value = 3179.8
$ws.Range("K22").Value = 3179.8
$ws.Range("M22").Value = -2884.8

$ws.Range("H27").Value = 3179.8
$ws.Range("I27").Value = 3179.8
$ws.Range("K27").Value = 3179.8
$ws.Range("M27").Value = -3072.8

$ws.Range("H46").Value = 15578.579
$ws.Range("I46").Value = 2814.8
$ws.Range("K46").Value = 2814.8
$ws.Range("M46").Value = -2626.8

$ws.Range("H53").Value = 30000
$ws.Range("I53").Value = 30000
$ws.Range("K53").Value = 30000
$ws.Range("M53").Value = -29482

$ws.Range("H82").Value = 13680.474
$ws.Range("I82").Value = 17539
$ws.Range("K82").Value = 17539
$ws.Range("M82").Value = -17178

$ws.Range("H85").Value = 13680.474
$ws.Range("I85").Value = 17539
$ws.Range("K85").Value = 17539
$ws.Range("M85").Value = -16291

$ws.Range("H100").Value = 12504438
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H132").Value = 19862
$ws.Range("I132").Value = 23316.777
$ws.Range("K132").Value = 69950.33099999999
$ws.Range("M132").Value = -67420.33099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 27600
$ws.Range("J54").Value = 27600
$ws.Range("L54").Value = 27600
$ws.Range("N54").Value = -28640

$ws.Range("H107").Value = 1154.8
$ws.Range("I107").Value = 1068.75
$ws.Range("K107").Value = 3206.25
$ws.Range("M107").Value = -1286.25

$ws.Range("H127").Value = 22033.334
$ws.Range("J127").Value = 22033.334
$ws.Range("L127").Value = 22033.334
$ws.Range("N127").Value = -31953.334

$ws.Range("H132").Value = 4980.4287
$ws.Range("I132").Value = 4594.3076
$ws.Range("K132").Value = 13782.9228
$ws.Range("M132").Value = -11252.9228
